$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "270.91"
    3  = "21.10"
    4  = "6.252"
    5  = "0.06199"
    6  = "3.559"
    7  = "6.540"
    8  = "1.432"
    9  = "0.8226"
    10 = "0.1651"
    11 = "0.08282"
    12 = "0.03553"
    13 = "0.03185"
    14 = "0.09184"
    15 = "3.769"
    16 = "0.001635"
    17 = "0.04669"
    18 = "0.006401"
    20 = "0.001068"
    23 = "2.301"
    24 = "0.01378"
    25 = "0.3290"
    28 = "0.0002712"
    40 = "0.04709"
    41 = "0.006997"
    42 = "0.004499"
    44 = "0.01048"
    45 = "0.00006238"
    46 = "0.0009903"
    48 = "0.9903"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
}
